$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to "2566"
$ws.Name = "2566"

# Delete column A (index numbers), shifting B->A, C->B, D->C
$ws.Columns.Item(1).Delete()

# Clear all cell formatting/styles across the used range (remove bold font, border, alignment)
$ws.Cells.ClearFormats()

# Remove the now-empty "Face Status" data cells (C2:C22), keeping only the header in C1
$ws.Range("C2:C22").Clear()
